# Apply cryptos list update (prices/volumes refreshed; three rows reordered)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value looks like a plain decimal number must be
# pre-formatted as Text, otherwise Excel auto-converts the literal into a real
# number (e.g. "297.59" -> 297.59) and the multi-dot/"NN.NN" text layout used
# throughout this sheet would be lost.
$textCells = $ws.Range("D5", "D6", "D7", "D9", "D10", "D12", "D16", "D17", "D20", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D30", "D31", "D32", "D33", "D34", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D45", "D47", "D48", "D49", "D50", "D51")
$textCells.NumberFormat = "@"

$ws.Range("D2").Value = "43.793.79"
$ws.Range("E2").Value = "  -1.45%  "

$ws.Range("D3").Value = "2.212.78"
$ws.Range("E3").Value = "  -1.18%  "

$ws.Range("E4").Value = "  -1.49%  "

$ws.Range("D5").Value = "297.59"
$ws.Range("E5").Value = "  -3.22%  "

$ws.Range("D6").Value = "89.46"
$ws.Range("E6").Value = "  -4.63%  "

$ws.Range("D7").Value = "0.554"
$ws.Range("E7").Value = "  -2.98%  "

$ws.Range("E8").Value = "  -0.60%  "

$ws.Range("D9").Value = "0.487"
$ws.Range("E9").Value = "  -6.61%  "

$ws.Range("D10").Value = "32.49"
$ws.Range("E10").Value = "  -6.09%  "

$ws.Range("E11").Value = "  -3.78%  "

$ws.Range("D12").Value = "6.90"
$ws.Range("E12").Value = "  -4.16%  "

$ws.Range("E13").Value = "  -0.93%  "

$ws.Range("D14").Value = "2.549.74"
$ws.Range("E14").Value = "  -1.14%  "

$ws.Range("D15").Value = "2.218.10"
$ws.Range("E15").Value = "  -4.14%  "

$ws.Range("D16").Value = "13.26"
$ws.Range("E16").Value = "  -1.94%  "

$ws.Range("D17").Value = "0.770"
$ws.Range("E17").Value = "  -7.41%  "

$ws.Range("D18").Value = "43.614.65"
$ws.Range("E18").Value = "  -1.03%  "

$ws.Range("D19").Value = "0.0₃0899"
$ws.Range("E19").Value = "  -5.71%  "

$ws.Range("D20").Value = "11.52"
$ws.Range("E20").Value = "  -2.45%  "

$ws.Range("E21").Value = "  -7.25%  "

$ws.Range("D22").Value = "64.19"
$ws.Range("E22").Value = "  -2.11%  "

$ws.Range("D23").Value = "234.84"
$ws.Range("E23").Value = "  -1.09%  "

$ws.Range("D24").Value = "2.79"
$ws.Range("E24").Value = "  -5.57%  "

$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.34%  "

$ws.Range("D26").Value = "1.84"
$ws.Range("E26").Value = "  -6.69%  "

$ws.Range("B27").Value = "InjectiveProtocol"
$ws.Range("C27").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D27").Value = "37.92"
$ws.Range("E27").Value = "  +0.37%  "

$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "2.19"
$ws.Range("E28").Value = "  -1.73%  "

$ws.Range("E29").Value = "  -4.91%  "

$ws.Range("D30").Value = "152.79"
$ws.Range("E30").Value = "  -0.43%  "

$ws.Range("D31").Value = "19.02"
$ws.Range("E31").Value = "  -4.59%  "

$ws.Range("D32").Value = "5.38"
$ws.Range("E32").Value = "  -9.12%  "

$ws.Range("D33").Value = "0.0752"
$ws.Range("E33").Value = "  -5.25%  "

$ws.Range("D34").Value = "2.47"
$ws.Range("E34").Value = "  -6.21%  "

$ws.Range("E35").Value = "  -2.30%  "

$ws.Range("D36").Value = "2.82"
$ws.Range("E36").Value = "  -9.00%  "

$ws.Range("D37").Value = "0.101"
$ws.Range("E37").Value = "  -8.67%  "

$ws.Range("D38").Value = "1.67"
$ws.Range("E38").Value = "  -7.78%  "

$ws.Range("D39").Value = "0.0297"
$ws.Range("E39").Value = "  -0.50%  "

$ws.Range("D40").Value = "3.57"
$ws.Range("E40").Value = "  -5.19%  "

$ws.Range("D41").Value = "3.12"
$ws.Range("E41").Value = "  -7.80%  "

$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").Value = "13.11"
$ws.Range("E42").Value = "  -11.90%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.93%  "

$ws.Range("D44").Value = "1.805.89"
$ws.Range("E44").Value = "  +1.94%  "

$ws.Range("D45").Value = "1.77"
$ws.Range("E45").Value = "  +12.89%  "

$ws.Range("E46").Value = "  -5.49%  "

$ws.Range("D47").Value = "94.07"
$ws.Range("E47").Value = "  -4.53%  "

$ws.Range("D48").Value = "66.30"
$ws.Range("E48").Value = "  -5.51%  "

$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D49").Value = "72.16"
$ws.Range("E49").Value = "  -8.42%  "

$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "7.70"
$ws.Range("E50").Value = "  -4.73%  "

$ws.Range("D51").Value = "4.55"
$ws.Range("E51").Value = "  -6.66%  "

# Drop back to the default style so formatting matches the source workbook
# (only the text content should differ from "before").
$textCells.Style = "Normal"
